# Applies the "Temporary_day" diet-plan update:
#  - Lunedi / Martedi / Mercoledi: header "Nome" -> "Alimento" (where applicable),
#    drop the trailing single food row, add a 4-row macro summary block.
#  - Venerdi / Sabato: populate the previously-empty sheets with a full food
#    list + macro summary.
#  - Temporary_day: brand new sheet appended at the end of the workbook with
#    its own food list + macro summary.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a "label -> value" summary block of 4 rows starting at
# $startRow on worksheet $ws (columns A/B, values already formatted strings).
# ---------------------------------------------------------------------------
function Write-Summary($ws, $startRow, $protein, $fat, $carbs, $kcal) {
    $ws.Range("A$startRow").Value = "Apporto proteico totale:"
    $ws.Range("B$startRow").Value = $protein
    $r = $startRow + 1
    $ws.Range("A$r").Value = "Apporto di grassi totale:"
    $ws.Range("B$r").Value = $fat
    $r = $startRow + 2
    $ws.Range("A$r").Value = "Apporto di carboidrati totale:"
    $ws.Range("B$r").Value = $carbs
    $r = $startRow + 3
    $ws.Range("A$r").Value = "Apporto di calorie totale:"
    $ws.Range("B$r").Value = $kcal
}

# ---------------------------------------------------------------------------
# Helper: fill header (Alimento / Valore (g)) + a list of [name, value] food
# rows starting at row 2, copying the header style from an existing sheet.
# ---------------------------------------------------------------------------
function Write-FoodSheet($ws, $foods) {
    $ws.Range("A1").Value = "Alimento"
    $ws.Range("B1").Value = "Valore (g)"
    $row = 2
    foreach ($food in $foods) {
        $ws.Range("A$row").Value = $food[0]
        $ws.Range("B$row").Value = $food[1]
        $row++
    }
}

# ===========================================================================
# Lunedi
# ===========================================================================
$ws = $wb.Worksheets.Item("Lunedi")
$ws.Range("A1").Value = "Alimento"
$ws.Range("A19:B19").ClearContents()
Write-Summary $ws 21 "228.4 g" "54.0 g" "208.4 g" "2300.0 kcal"

# ===========================================================================
# Martedi
# ===========================================================================
$ws = $wb.Worksheets.Item("Martedi")
$ws.Range("A1").Value = "Alimento"
$ws.Range("A18:B18").ClearContents()
Write-Summary $ws 20 "227.9 g" "54.0 g" "214.5 g" "2300.0 kcal"

# ===========================================================================
# Mercoledi
# ===========================================================================
$ws = $wb.Worksheets.Item("Mercoledi")
$ws.Range("A19:C22").ClearContents()
Write-Summary $ws 20 "217.5 g" "46.1 g" "229.4 g" "2300.0 kcal"

# ===========================================================================
# Venerdi (was empty)
# ===========================================================================
$ws = $wb.Worksheets.Item("Venerdi")
$lun = $wb.Worksheets.Item("Lunedi")
$lun.Range("A1:B1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

$foods = @(
    ,@("Milk pro budino (20g)", 200)
    ,@("banana", 120)
    ,@("bresaola", 70)
    ,@("cavolfiore", 100)
    ,@("fiocchi latte (conad)", 100)
    ,@("gocce di cioccolato", 5)
    ,@("insalata", 100)
    ,@("miele", 10)
    ,@("muesli conad", 40)
    ,@("olio", 1)
    ,@("pasta integrale", 100)
    ,@("petto di pollo", 300)
    ,@("proteine buone", 30)
    ,@("ricotta di mucca", 300)
    ,@("riso", 58.1)
    ,@("wasa", 13)
)
Write-FoodSheet $ws $foods
Write-Summary $ws 20 "219.0 g" "54.0 g" "220.2 g" "2300.0 kcal"

# ===========================================================================
# Sabato (was empty)
# ===========================================================================
$ws = $wb.Worksheets.Item("Sabato")
$lun.Range("A1:B1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

$foods = @(
    ,@("Milk pro budino (20g)", 200)
    ,@("albume", 100)
    ,@("barilla lenticchie ", 100)
    ,@("bresaola", 60)
    ,@("cavolfiore", 100)
    ,@("farina (normale)", 100)
    ,@("fettine vitello", 300)
    ,@("insalata", 100)
    ,@("olio", 20)
    ,@("proteine buone", 30)
    ,@("riso", 50)
    ,@("sciroppo acero", 25.4)
    ,@("stracchino", 160)
    ,@("uova", 0)
    ,@("wasa", 13)
    ,@("yogurt 0% bianco", 100)
)
Write-FoodSheet $ws $foods
Write-Summary $ws 19 "216.6 g" "51.8 g" "220.5 g" "2300.0 kcal"

# ===========================================================================
# Temporary_day (brand new sheet, appended at the very end)
# ===========================================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Temporary_day"
$lun.Range("A1:B1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

$foods = @(
    ,@("Milk pro budino (20g)", 200)
    ,@("barretta proteica", 64)
    ,@("bresaola", 60)
    ,@("cavolfiore", 100)
    ,@("fiocchi avena", 100)
    ,@("gocce di cioccolato", 10)
    ,@("insalata", 100)
    ,@("latte granarolo", 100)
    ,@("olio", 7.7)
    ,@("petto di pollo", 250)
    ,@("riso", 122.7)
    ,@("salmone bollito", 250)
    ,@("wasa", 13)
)
Write-FoodSheet $ws $foods
$ws.Range("A16").Value = "yogurt 0% bianco"
$ws.Range("B16").Value = 100
Write-Summary $ws 17 "216.5 g" "54.0 g" "219.7 g" "2300.0 kcal"
